$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the duplicated bold heading paragraph that was accidentally left
#    at the bottom of the document ("Play Chicken Fox Free - Exciting
#    Multipliers and Free Spins!") and rewrite the final (italic) paragraph
#    with the new image-prompt text, keeping its original formatting
#    (leading empty run + italic run).
# ---------------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastText = $lastPara.Range.Text

if ($lastText -like "*Experience the thrill of big wins*") {
    $bottomBold = $d.Paragraphs.Item($lastIndex - 1)
    if ($bottomBold.Range.Text -like "*Play Chicken Fox Free*") {
        $bottomBold.Range.Delete()
    }

    $finalIndex = $d.Paragraphs.Count
    $finalPara = $d.Paragraphs.Item($finalIndex)
    $rStart = $finalPara.Range.Start
    $rEnd = $finalPara.Range.End
    $target = $d.Range($rStart, $rEnd - 1)
    $target.Text = "Create a feature image for Chicken Fox that captures the cute farmyard theme of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be depicted alongside the game's farm animals, such as chickens, goats, and pigs, with a background of a traditional farm setting, complete with a barn and rolling hills. The image should convey the fun, playful nature of the game, while also highlighting its unique features, such as the Free Games feature and multipliers."
}

# ---------------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$newRange = $titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the thrill of big wins! Play Chicken Fox, an online slot game with multipliers, wilds, and free spins. Try it for free now!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $metaPara.Range.InsertXML($metaXml)
